$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells A1:D1 to clean snake_case column names
$cellUpdates = @(
    @("A1", "mx_state"),
    @("B1", "mx_municipality"),
    @("C1", "n_matriculas"),
    @("D1", "pct_matriculas"),
    @("B6", "Pabellón De Arteaga"),
    @("B7", "Rincón De Romos"),
    @("B8", "San Francisco De Los Romo"),
    @("B9", "San José De Gracia"),
    @("B25", "Amatenango De La Frontera"),
    @("B34", "Comitán De Domínguez"),
    @("B54", "Mazapa De Madero"),
    @("B63", "San Cristóbal De Las Casas"),
    @("B96", "Coyame Del Sotol"),
    @("B105", "Guadalupe Y Calvo"),
    @("B108", "Hidalgo Del Parral"),
    @("B129", "San Francisco De Borja"),
    @("B130", "San Francisco Del Oro"),
    @("B138", "Valle De Zaragoza"),
    @("B150", "San Juan De Sabinas"),
    @("B161", "Villa De Álvarez"),
    @("A163", "Ciudad De México"),
    @("B167", "Cuajimalpa De Morelos"),
    @("B182", "Coneto De Comonfort"),
    @("B195", "Nombre De Dios"),
    @("B199", "Pánuco De Coronado"),
    @("B206", "San Juan De Guadalupe"),
    @("B207", "San Juan Del Río"),
    @("B208", "San Pedro Del Gallo"),
    @("A218", "Estado De México"),
    @("B218", "Acambay De Ruíz Castañeda"),
    @("B221", "Almoloya De Alquisiras"),
    @("B222", "Almoloya De Juárez"),
    @("B226", "Atizapán De Zaragoza"),
    @("B231", "Chapa De Mota"),
    @("B241", "Ecatepec De Morelos"),
    @("B244", "Ixtapan De La Sal"),
    @("B245", "Ixtapan Del Oro"),
    @("B257", "Naucalpan De Juárez"),
    @("B265", "San Felipe Del Progreso"),
    @("B273", "Tenango Del Valle"),
    @("B281", "Tlalnepantla De Baz"),
    @("B286", "Valle De Bravo"),
    @("B287", "Valle De Chalco Solidaridad"),
    @("B288", "Villa De Allende"),
    @("B289", "Villa Del Carbón"),
    @("B301", "San Miguel De Allende"),
    @("B302", "Apaseo El Alto"),
    @("B303", "Apaseo El Grande"),
    @("B311", "Dolores Hidalgo Cuna De La Independencia Nacional"),
    @("B315", "Jaral Del Progreso"),
    @("B323", "Purísima Del Rincón"),
    @("B328", "San Francisco Del Rincón"),
    @("B330", "San Luis De La Paz"),
    @("B332", "Silao De La Victoria"),
    @("B336", "Valle De Santiago"),
    @("B341", "Acapulco De Juárez"),
    @("B343", "Ajuchitlán Del Progreso"),
    @("B344", "Alcozauca De Guerrero"),
    @("B348", "Atenango Del Río"),
    @("B350", "Atoyac De Álvarez"),
    @("B351", "Ayutla De Los Libres"),
    @("B354", "Buenavista De Cuéllar"),
    @("B355", "Chilapa De Álvarez"),
    @("B356", "Chilpancingo De Los Bravo"),
    @("B357", "Coahuayutla De José María Izazaga"),
    @("B362", "Coyuca De Benítez"),
    @("B363", "Coyuca De Catalán"),
    @("B367", "Cuetzala Del Progreso"),
    @("B368", "Cutzamala De Pinzón"),
    @("B374", "Huitzuco De Los Figueroa"),
    @("B375", "Iguala De La Independencia"),
    @("B377", "Ixcateopan De Cuauhtémoc"),
    @("B378", "Zihuatanejo De Azueta"),
    @("B380", "La Unión De Isidoro Montes De Oca"),
    @("B395", "Taxco De Alarcón"),
    @("B397", "Técpan De Galeana"),
    @("B399", "Tepecoacuilco De Trujano"),
    @("B401", "Tixtla De Guerrero"),
    @("B405", "Tlalixtaquilla De Maldonado"),
    @("B406", "Tlapa De Comonfort"),
    @("B414", "Agua Blanca De Iturbide"),
    @("B418", "Atotonilco De Tula"),
    @("B419", "Atotonilco El Grande"),
    @("B423", "Cuautepec De Hinojosa"),
    @("B426", "Huasca De Ocampo"),
    @("B429", "Jacala De Ledezma"),
    @("B434", "Mineral De La Reforma"),
    @("B435", "Mineral Del Monte"),
    @("B436", "Mixquiahuala De Juárez"),
    @("B437", "Nopala De Villagrán"),
    @("B438", "Omitlán De Juárez"),
    @("B439", "Pachuca De Soto"),
    @("B440", "Progreso De Obregón"),
    @("B445", "Santiago Tulantepec De Lugo Guerrero"),
    @("B449", "Tenango De Doria"),
    @("B451", "Tepeji Del Río De Ocampo"),
    @("B453", "Tezontepec De Aldama"),
    @("B457", "Tula De Allende"),
    @("B458", "Tulancingo De Bravo"),
    @("B459", "Zacualtipán De Ángeles"),
    @("B460", "Zapotlán De Juárez"),
    @("B464", "Acatlán De Juárez"),
    @("B469", "Atotonilco El Alto"),
    @("B470", "Autlán De Navarro"),
    @("B475", "Cañadas De Obregón"),
    @("B481", "Cuautitlán De García Barragán"),
    @("B489", "Encarnación De Díaz"),
    @("B494", "Huejuquilla El Alto"),
    @("B495", "Ixtlahuacán Del Río"),
    @("B498", "Jilotlán De Los Dolores"),
    @("B504", "Lagos De Moreno"),
    @("B510", "Ojuelos De Jalisco"),
    @("B514", "San Cristóbal De La Barranca"),
    @("B515", "San Diego De Alejandría"),
    @("B517", "San Juan De Los Lagos"),
    @("B521", "San Miguel El Alto"),
    @("B522", "Santa María De Los Ángeles"),
    @("B523", "Santa María Del Oro"),
    @("B526", "Talpa De Allende"),
    @("B527", "Tamazula De Gordiano"),
    @("B531", "Teocuitatlán De Corona"),
    @("B532", "Tepatitlán De Morelos"),
    @("B534", "Tizapán El Alto"),
    @("B535", "Tlajomulco De Zúñiga"),
    @("B546", "Unión De San Antonio"),
    @("B547", "Unión De Tula"),
    @("B551", "Yahualica De González Gallo"),
    @("B552", "Zacoalco De Torres"),
    @("B555", "Zapotitlán De Vadillo"),
    @("B556", "Zapotlán Del Rey"),
    @("B557", "Zapotlán El Grande"),
    @("B575", "Coalcomán De Vázquez Pallares"),
    @("B577", "Cojumatlán De Régules"),
    @("B667", "Puente De Ixtla"),
    @("B672", "Tlaltizapán De Zapata"),
    @("B678", "Zacualpan De Amilpas"),
    @("B682", "Amatlán De Cañas"),
    @("B683", "Bahía De Banderas"),
    @("B687", "Ixtlán Del Río"),
    @("B694", "Santa María Del Oro"),
    @("B712", "San Nicolás De Los Garza"),
    @("B716", "Ayoquezco De Aldama"),
    @("B718", "Chalcatongo De Hidalgo"),
    @("B719", "Ciénega De Zimatlán"),
    @("B720", "Coicoyán De Las Flores"),
    @("B723", "Cuilápam De Guerrero"),
    @("B724", "Heroica Ciudad De Ejutla De Crespo"),
    @("B725", "Heroica Ciudad De Huajuapan De León"),
    @("B726", "Heroica Ciudad De Tlaxiaco"),
    @("B727", "Huautla De Jiménez"),
    @("B728", "Ixtlán De Juárez"),
    @("B729", "Heroica Ciudad De Juchitán De Zaragoza"),
    @("B735", "Mariscala De Juárez"),
    @("B738", "Oaxaca De Juárez"),
    @("B739", "Ocotlán De Morelos"),
    @("B740", "Pinotepa De Don Luis"),
    @("B741", "Putla Villa De Guerrero"),
    @("B780", "San Miguel Del Puerto"),
    @("B793", "San Pedro Y San Pablo Tequixtepec"),
    @("B807", "Santa María Del Tule"),
    @("B810", "Santa María Jalapa Del Marqués"),
    @("B833", "Teotitlán De Flores Magón"),
    @("B834", "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"),
    @("B835", "Tlacolula De Matamoros"),
    @("B836", "Villa De Etla"),
    @("B837", "Villa De Tututepec"),
    @("B838", "Villa De Zaachila"),
    @("B839", "Villa Sola De Vega"),
    @("B840", "Zimatlán De Álvarez"),
    @("B855", "Chalchicomula De Sesma"),
    @("B861", "Cuayuca De Andrade"),
    @("B870", "Ixcamilpa De Guerrero"),
    @("B872", "Izúcar De Matamoros"),
    @("B878", "Los Reyes De Juárez"),
    @("B880", "Palmar De Bravo"),
    @("B893", "San Nicolás De Los Ranchos"),
    @("B895", "San Salvador El Seco"),
    @("B904", "Tepexi De Rodríguez"),
    @("B906", "Tetela De Ocampo"),
    @("B911", "Tlacotepec De Benito Juárez"),
    @("B918", "Totoltepec De Guerrero"),
    @("B929", "Amealco De Bonfil"),
    @("B930", "Cadereyta De Montes"),
    @("B934", "Landa De Matamoros"),
    @("B937", "Pinal De Amoles"),
    @("B940", "San Juan Del Río"),
    @("B946", "Armadillo De Los Infante"),
    @("B954", "San Ciro De Acosta"),
    @("B957", "Santa María Del Río"),
    @("B959", "Soledad De Graciano Sánchez"),
    @("B963", "Villa De Arriaga"),
    @("B964", "Villa De Ramos"),
    @("B965", "Villa De Reyes"),
    @("B1030", "Nacozari De García"),
    @("B1041", "San Felipe De Jesús"),
    @("B1044", "San Pedro De La Cueva"),
    @("B1063", "Jalpa De Méndez"),
    @("B1078", "Soto La Marina"),
    @("B1085", "Contla De Juan Cuamatzi"),
    @("B1088", "Ixtacuixtla De Mariano Matamoros"),
    @("B1093", "Papalotla De Xicohténcatl"),
    @("B1102", "Alto Lucero De Gutiérrez Barrios"),
    @("B1105", "Amatlán De Los Reyes"),
    @("B1111", "Boca Del Río"),
    @("B1115", "Castillo De Teayo"),
    @("B1117", "Cazones De Herrera"),
    @("B1124", "Cosamaloapan De Carpio"),
    @("B1137", "Hueyapan De Ocampo"),
    @("B1139", "Ixhuatlán De Madero"),
    @("B1140", "Ixhuatlán Del Sureste"),
    @("B1148", "Lerdo De Tejada"),
    @("B1150", "Martínez De La Torre"),
    @("B1151", "Medellín De Bravo"),
    @("B1160", "Paso De Ovejas"),
    @("B1163", "Poza Rica De Hidalgo"),
    @("B1170", "Sayula De Alemán"),
    @("B1172", "Soledad De Doblado"),
    @("B1182", "Tlacotepec De Mejía"),
    @("B1187", "Vega De Alatorre"),
    @("B1205", "Cañitas De Felipe Pescador"),
    @("B1207", "Concepción Del Oro"),
    @("B1209", "El Plateado De Joaquín Amaro"),
    @("B1218", "Jiménez Del Teul"),
    @("B1227", "Moyahua De Estrada"),
    @("B1228", "Nochistlán De Mejía"),
    @("B1237", "Teúl De González Ortega"),
    @("B1238", "Tlaltenango De Sánchez Román"),
    @("B1240", "Villa De Cos")
)

foreach ($pair in $cellUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Remove the trailing footnote/metadata rows (now rows 1248-1252)
$ws.Rows("1248:1252").Delete()

